# "Generate Report for Handoff"
# Updates the handoff/generation timestamps produced for the last row
# (7242f80f-4118-4de1-b9a1-cc7307253753.md) on each localization sheet,
# plus the aggregated "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

# zh-cn: "Latest Handoff Datetime" (column H) for row 7 gets a fresh handoff timestamp.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-05 06:47:11"

# de-de: "Latest Handoff Datetime" (column H) for row 7 gets a fresh handoff timestamp.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-05 06:47:16"

# Overview: "Latest HO Xliff Generate Date" (column G) for row 7 reflects the newest
# handoff generation across languages.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-05 06:47:16"
